# Children.xlsx — update "Apgar" header split + NPR column relabeling,
# resize several columns, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row relabeling (row 1) ---
# Apgar score is now split into a 1-minute and a 5-minute reading.
$ws.Range("F1").Value = "Апгар(1 мин)"
$ws.Range("G1").Value = "Апгар(5 мин)"

# The NPR sub-columns are now explicitly "delay" (Задержка) columns.
$ws.Range("O1").Value = "Задержка моторика"
$ws.Range("P1").Value = "Задержка речь"
$ws.Range("Q1").Value = "Задержка моторика + речь"

# Physical development column renamed (set last so the new shared string
# is appended after the three delay columns above).
$ws.Range("M1").Value = "Норма физ. Развития"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 5.33   # A
$ws.Columns.Item(2).ColumnWidth = 14.5   # B
$ws.Columns.Item(3).ColumnWidth = 6      # C
$ws.Columns.Item(6).ColumnWidth = 10.67  # F
$ws.Columns.Item(7).ColumnWidth = 11     # G
$ws.Columns.Item(15).ColumnWidth = 17.67 # O
$ws.Columns.Item(16).ColumnWidth = 13.33 # P
$ws.Columns.Item(17).ColumnWidth = 23.67 # Q

# --- Selection moves from G2 to L27 ---
$ws.Range("L27").Select()
